$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 82: new "2023" year-header label in column A, reusing the same look as the
# existing year header in A68 ("2022"). Use a leading apostrophe so the numeric-
# looking text "2023" is stored as text (shared string), not as a number.
$ws.Range("A68").Copy()
$ws.Range("A82").PasteSpecial(-4122)
$ws.Range("A82").Value = "'2023"

# Row 83: January 2023 entry (date, particulars, earned days)
$ws.Range("A83").Value = 44927
$ws.Range("B83").Value = "SP(1-0-00)"
$ws.Range("C83").Value = 1.25

# K83 remarks date - reuse the style already used by the other remark dates
# (e.g. K71) so the same cell style index is applied instead of a brand new one.
$ws.Range("K71").Copy()
$ws.Range("K83").PasteSpecial(-4122)
$ws.Range("K83").Value = 44937

# Row 84: February 2023 entry
$ws.Range("A84").Value = 44958
$ws.Range("C84").Value = 1.25

# Rows 85-120: continue the month-by-month date sequence in column A only
# (March 2023 through February 2026), one row per month, 1st of month.
$dateSerials = @(
    44986, 45017, 45047, 45078, 45108, 45139, 45170, 45200, 45231, 45261,
    45292, 45323, 45352, 45383, 45413, 45444, 45474, 45505, 45536, 45566,
    45597, 45627, 45658, 45689, 45717, 45748, 45778, 45809, 45839, 45870,
    45901, 45931, 45962, 45992, 46023, 46054
)
for ($i = 0; $i -lt $dateSerials.Length; $i++) {
    $row = 85 + $i
    $ws.Range("A$row").Value = $dateSerials[$i]
}

$excel.CutCopyMode = $false

# Scroll the bottom (frozen/split) pane down so row 79 becomes its top row, and
# leave the selection on B84 - matching where the user was last working.
$excel.ActiveWindow.SplitRow = 78
$ws.Range("B84").Select() | Out-Null
